$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns("D").Insert()
$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 43404
